# Kraftwerksgrafik in Visualisierung eingebaut
#
# Applies the changes described by the diff to Kraftwerke_Lasten_Speichertabelle.xlsx:
#  - L7: 70 -> 0.7
#  - L8: 100 -> 1
#  - Add a new data row 14 (power plant #13) with values in A14:O14
#  - Extend the shared formulas in columns H and I down to row 14
#  - Update the active selection to P14 (and scroll the view toward column H)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct variable cost entries (were entered as percent*100 instead of fraction) ---
$ws.Range("L7").Value = 0.7
$ws.Range("L8").Value = 1

# --- New power plant entry in row 14 ---
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 7
$ws.Range("C14").Value = 100
$ws.Range("D14").Value = -1
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 1

# Columns H and I carry shared formulas (4*C.. and 5*F..) down through row 13;
# extend that same calculation down into the new row 14.
$ws.Range("H14").Formula = "=4*C14"
$ws.Range("I14").Formula = "=5*F14"

$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("O14").Value = 0

# --- Update the window/selection state to reflect where the user ended up working ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 8
$win.ScrollRow = 1
$ws.Range("P14").Select() | Out-Null
